$d = $word.ActiveDocument

# The only real content change in this revision is the generation
# timestamp stamped into the document footer ("PubMedOut-4" is a freshly
# regenerated export, so its footer banner moves from 01:36Z to 12:13Z).
foreach ($sec in $d.Sections) {
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $ftr.Range.Find.Execute(
                "2025-06-30 01:36Z / ", $true, $false, $false, $false, $false,
                $true, 1, $false, "2025-06-30 12:13Z / ", 2) | Out-Null
        }
    }
}
